$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '41.669.07'
$ws.Range("E2").Value = '  +0.34%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.476.72'
$ws.Range("E3").Value = '  +0.08%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '317.53'
$ws.Range("E5").Value = '  +1.68%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '92.43'
$ws.Range("E6").Value = '  +0.59%  '
$ws.Range("E7").Value = '  +0.74%  '
$ws.Range("E8").Value = '  +0.06%  '
$ws.Range("E9").Value = '  +0.52%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '33.14'
$ws.Range("E10").Value = '  +1.84%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0857'
$ws.Range("E11").Value = '  +8.56%  '
$ws.Range("E12").Value = '  +0.35%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '2.857.33'
$ws.Range("E13").Value = '  -0.08%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.91'
$ws.Range("E14").Value = '  +0.52%  '
$ws.Range("E15").Value = '  -4.23%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.479.17'
$ws.Range("E16").Value = '  -0.45%  '
$ws.Range("E17").Value = '  +2.60%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '41.617.56'
$ws.Range("E18").Value = '  +0.27%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.48'
$ws.Range("E19").Value = '  -0.52%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0950'
$ws.Range("E20").Value = '  +1.16%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '71.39'
$ws.Range("E21").Value = '  -0.83%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '11.34'
$ws.Range("E22").Value = '  +2.64%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '240.10'
$ws.Range("E23").Value = '  +1.65%  '
$ws.Range("E24").Value = '  +0.91%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.94'
$ws.Range("E25").Value = '  +1.97%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.999'
$ws.Range("E26").Value = '  -0.08%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '24.73'
$ws.Range("E27").Value = '  -1.04%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.28'
$ws.Range("E28").Value = '  +2.68%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.88'
$ws.Range("E29").Value = '  +2.11%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '36.37'
$ws.Range("E30").Value = '  +1.91%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '160.19'
$ws.Range("E31").Value = '  +1.91%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.52'
$ws.Range("E32").Value = '  +1.61%  '
$ws.Range("E33").Value = '  +0.05%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.58'
$ws.Range("E34").Value = '  +0.36%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0767'
$ws.Range("E35").Value = '  +1.34%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '17.43'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.87'
$ws.Range("E37").Value = '  +3.06%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.94'
$ws.Range("E38").Value = '  +1.49%  '
$ws.Range("E39").Value = '  +1.91%  '
$ws.Range("E40").Value = '  -1.13%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.99'
$ws.Range("E41").Value = '  -1.10%  '
$ws.Range("E42").Value = '  +2.76%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.988.01'
$ws.Range("E43").Value = '  +1.34%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0286'
$ws.Range("E44").Value = '  +0.55%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '18.94'
$ws.Range("E45").Value = '  +0.02%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.00'
$ws.Range("E46").Value = '  +2.17%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.26'
$ws.Range("E47").Value = '  +3.32%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.713.70'
$ws.Range("E48").Value = '  -0.20%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '97.55'
$ws.Range("E49").Value = '  +0.01%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '74.29'
$ws.Range("E50").Value = '  +2.99%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '67.39'
$ws.Range("E51").Value = '  -0.22%  '
